$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1
#    title paragraph ("Play All That Cash Slot for Free - Review &
#    Ratings").
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs(2).Range

$metaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review on All That Cash, a unique 3 payline slot game. Learn about the Power Bet feature and RTP percentage. Play All That Cash slot for free.</w:t></w:r></w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
$metaPara.InsertXML($metaXml)

# ------------------------------------------------------------------
# 2) Near the end of the document, the paragraph that duplicated the
#    bold title ("Play All That Cash Slot for Free - Review &
#    Ratings") is removed, and the following italic paragraph's text
#    is replaced with the new feature-image prompt (keeping italics).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldPara = $null
$italicPara = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Play All That Cash Slot for Free - Review & Ratings" -and $i -gt 1) {
        $boldPara = $p
        $italicPara = $d.Paragraphs($i + 1)
        break
    }
}

$rangeStart = $boldPara.Range.Start
$rangeEnd = $italicPara.Range.End
$fullRange = $d.Range($rangeStart, $rangeEnd)

$newImagePrompt = "Create a feature image for All That Cash that features a happy Maya warrior with glasses in a cartoon style. The image should be bright and colorful with a green color scheme to represent money and wealth. The Maya warrior should be standing confidently with a big smile on their face and wearing glasses to represent intelligence and success. The background of the image can feature floating dollar bills and other money-themed elements to reinforce the theme of the game. The image should capture the excitement and fun of playing a slot game, while also highlighting the potential for big wins and the unique Power Bet feature of the game."

$endXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>$newImagePrompt</w:t></w:r></w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@
$fullRange.InsertXML($endXml)
